$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L2").Value = "[0.3691006962088813, 0.4699114309755926]"
$ws.Range("M2").Value = [double]"6.616929226765933e-14"
$ws.Range("N2").Value = [double]"6.616929226765933e-14"
$ws.Range("T2").Value = "[0.4840495674418006, 0.5362079588794344]"
